# Generate Report for Handoff
#
# Mirrors a localization-status report regeneration:
#   - marks rows 7,10,11,12,13,14 (on the zh-cn and de-de sheets) with
#     Priority = "ht" (handoff type)
#   - refreshes the "Latest Handoff Datetime" timestamp for those same rows
#   - refreshes the matching "Latest HO Xliff Generate Date" on the Overview
#     sheet for the one source file represented there (row 7)

$wb = $excel.ActiveWorkbook

$rows = @(7, 10, 11, 12, 13, 14)

# --- zh-cn sheet -----------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $wsZh.Range("E$r").Value = "ht"
    $wsZh.Range("H$r").Value = "2016-09-03 08:23:51"
}

# --- de-de sheet -------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $wsDe.Range("E$r").Value = "ht"
    $wsDe.Range("H$r").Value = "2016-09-03 08:23:55"
}

# --- Overview sheet ------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G7").Value = "2016-09-03 08:23:55"
